$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.779.02'
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").Value = '3.152.38'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '576.98'
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("D6").Value = '148.84'
$ws.Range("E6").Value = '  -1.79%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '3.152.51'
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("D9").Value = '0.525'
$ws.Range("E9").Value = '  -0.43%  '
$ws.Range("E10").Value = '  -2.02%  '
$ws.Range("E11").Value = '  -1.31%  '
$ws.Range("E12").Value = '  -0.87%  '
$ws.Range("D13").Value = '0.0000261'
$ws.Range("E13").Value = '  +2.73%  '
$ws.Range("D14").Value = '37.06'
$ws.Range("E14").Value = '  -1.29%  '
$ws.Range("D15").Value = '3.668.03'
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("D16").Value = '64.919.91'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Value = '3.153.02'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").Value = '7.11'
$ws.Range("E18").Value = '  -1.63%  '
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("D20").Value = '504.04'
$ws.Range("E20").Value = '  -1.68%  '
$ws.Range("D21").Value = '14.84'
$ws.Range("E21").Value = '  -0.69%  '
$ws.Range("D22").Value = '15.31'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '0.713'
$ws.Range("E23").Value = '  -3.32%  '
$ws.Range("D24").Value = '7.70'
$ws.Range("E24").Value = '  -1.60%  '
$ws.Range("D25").Value = '83.94'
$ws.Range("E25").Value = '  -1.32%  '
$ws.Range("D26").Value = '0.996'
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("D27").Value = '8.90'
$ws.Range("E27").Value = '  +2.10%  '
$ws.Range("E28").Value = '  -1.04%  '
$ws.Range("E29").Value = '  -1.02%  '
$ws.Range("D30").Value = '2.82'
$ws.Range("E30").Value = '  +6.26%  '
$ws.Range("D31").Value = '27.46'
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("E33").Value = '  +0.80%  '
$ws.Range("D34").Value = '6.15'
$ws.Range("E34").Value = '  +1.17%  '
$ws.Range("D35").Value = '6.45'
$ws.Range("E35").Value = '  -1.89%  '
$ws.Range("D36").Value = '54.51'
$ws.Range("E36").Value = '  -1.72%  '
$ws.Range("D37").Value = '0.0893'
$ws.Range("E37").Value = '  +3.65%  '
$ws.Range("D38").Value = '475.31'
$ws.Range("E38").Value = '  -1.81%  '
$ws.Range("D39").Value = '0.0415'
$ws.Range("E39").Value = '  -2.16%  '
$ws.Range("D40").Value = '2.93'
$ws.Range("E40").Value = '  -3.12%  '
$ws.Range("D41").Value = '8.67'
$ws.Range("E41").Value = '  +0.34%  '
$ws.Range("D42").Value = '3.005.97'
$ws.Range("E42").Value = '  -3.66%  '
$ws.Range("E43").Value = '  -3.87%  '
$ws.Range("D44").Value = '2.42'
$ws.Range("E44").Value = '  -1.45%  '
$ws.Range("D45").Value = '0.281'
$ws.Range("E45").Value = '  -3.54%  '
$ws.Range("D46").Value = '28.09'
$ws.Range("E46").Value = '  -3.90%  '
$ws.Range("E47").Value = '  +1.18%  '
$ws.Range("D50").Value = '34.56'
$ws.Range("E50").Value = '  +10.39%  '
$ws.Range("E51").Value = '  -2.94%  '
